$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A cells to Text format so numeric-looking strings stay as text (shared strings)
$ws.Range("A8:A10").NumberFormat = "@"

# Add rows 8, 9, 10 with channel labels (text) in column A and counts in column B
$ws.Range("A8").Value = "7008"
$ws.Range("B8").Value = 0

$ws.Range("A9").Value = "7009"
$ws.Range("B9").Value = 0

$ws.Range("A10").Value = "7010"
$ws.Range("B10").Value = 1

# Restore default General format/style so cells match the original style (s="0")
$ws.Range("A8:A10").NumberFormat = "General"
$ws.Range("A8:A10").Style = "Normal"
